$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("F6").Value = 7751
$ws.Range("F9").Value = 1147
$ws.Range("F12").Value = 23
$ws.Range("F13").Value = 173
$ws.Range("F17").Value = 2374
$ws.Range("F18").Value = 1022
$ws.Range("F22").Value = 6355
$ws.Range("F23").Value = 7003
$ws.Range("F24").Value = 402
$ws.Range("F32").Value = 519
$ws.Range("F33").Value = 519
$ws.Range("F39").Value = 419
$ws.Range("F41").Value = 1261
$ws.Range("F42").Value = 3264
$ws.Range("F46").Value = 46
$ws.Range("F48").Value = 6
$ws = $wb.Worksheets.Item(2)
$ws.Range("F7").Value = 88
$ws.Range("F10").Value = 294
$ws.Range("F11").Value = 3
$ws.Range("F23").Value = 4
$ws.Range("F25").Value = 6
$ws.Range("F26").Value = 6617
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 2004
$ws.Range("F5").Value = 1325
$ws.Range("F7").Value = 559
$ws.Range("F8").Value = 2152
$ws.Range("F9").Value = 8941
$ws.Range("F10").Value = 1071
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 2004
$ws.Range("F5").Value = 7751
$ws.Range("F6").Value = 1325
$ws.Range("F8").Value = 1071
$ws.Range("F10").Value = 1147
$ws.Range("F13").Value = 88
$ws.Range("F14").Value = 173
$ws.Range("F16").Value = 3
$ws.Range("F18").Value = 2374
$ws.Range("F19").Value = 1022
$ws.Range("F22").Value = 6355
$ws.Range("F23").Value = 7003
$ws.Range("F24").Value = 402
$ws.Range("F31").Value = 519
$ws.Range("F37").Value = 419
$ws.Range("F40").Value = 3264
$ws.Range("F45").Value = 6
